$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 288-289 (existing rows 288-309 shift down to 290-311).
$ws.Rows("288:289").Insert()

# New record (Primera quality) reported for this period.
$ws.Cells.Item(288, 1).Value = 7
$ws.Cells.Item(288, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(288, 3).Value = "Ñuble"
$ws.Cells.Item(288, 4).Value = 45013
$ws.Cells.Item(288, 5).Value = 16
$ws.Cells.Item(288, 6).Value = 100112017
$ws.Cells.Item(288, 7).Value = "Apio"
$ws.Cells.Item(288, 8).Value = "Americana (o)"
$ws.Cells.Item(288, 9).Value = "Primera"
$ws.Cells.Item(288, 10).Value = 40
$ws.Cells.Item(288, 11).Value = 7000
$ws.Cells.Item(288, 12).Value = 7000
$ws.Cells.Item(288, 13).Value = 7000
$ws.Cells.Item(288, 14).Value = "$/docena de matas"
$ws.Cells.Item(288, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(288, 16).Value = 1167
$ws.Cells.Item(288, 17).Value = 6
$ws.Cells.Item(288, 18).Value = "Hortaliza"

# New record (Segunda quality) reported for this period.
$ws.Cells.Item(289, 1).Value = 7
$ws.Cells.Item(289, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(289, 3).Value = "Ñuble"
$ws.Cells.Item(289, 4).Value = 45013
$ws.Cells.Item(289, 5).Value = 16
$ws.Cells.Item(289, 6).Value = 100112017
$ws.Cells.Item(289, 7).Value = "Apio"
$ws.Cells.Item(289, 8).Value = "Americana (o)"
$ws.Cells.Item(289, 9).Value = "Segunda"
$ws.Cells.Item(289, 10).Value = 20
$ws.Cells.Item(289, 11).Value = 6000
$ws.Cells.Item(289, 12).Value = 6000
$ws.Cells.Item(289, 13).Value = 6000
$ws.Cells.Item(289, 14).Value = "$/docena de matas"
$ws.Cells.Item(289, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(289, 16).Value = 1000
$ws.Cells.Item(289, 17).Value = 6
$ws.Cells.Item(289, 18).Value = "Hortaliza"

# Two existing records (now at rows 297 and 299 after the shift) had their "Origen"
# corrected from "Región de Coquimbo" to "Provincia del Elquí".
$ws.Cells.Item(297, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(299, 15).Value = "Provincia del Elquí"
